$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking price cells as Text so they round-trip as literal strings
# (matches how the source data was originally authored as inline text, not numbers).
$textCells = 'D5', 'D6', 'D7', 'D8', 'D13', 'D16', 'D20', 'D21', 'D24', 'D28', 'D31', 'D32', 'D35', 'D38', 'D40', 'D41', 'D43', 'D45', 'D46', 'D47'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.542.94'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '2.605.82'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '514.97'
$ws.Range('E5').Value = '  +1.99%  '
$ws.Range('D6').Value = '153.83'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  +4.66%  '
$ws.Range('D9').Value = '2.617.74'
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('E10').Value = '  +3.51%  '
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').Value = '0.130'
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').Value = '3.062.15'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').Value = '60.539.65'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '21.66'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').Value = '2.612.71'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = '358.00'
$ws.Range('E20').Value = '  +5.07%  '
$ws.Range('D21').Value = '10.61'
$ws.Range('E21').Value = '  +2.16%  '
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '61.08'
$ws.Range('E24').Value = '  +2.01%  '
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('D26').Value = '2.725.37'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('D28').Value = '0.995'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').Value = '0.0₃0841'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('E30').Value = '  -1.64%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '19.44'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('E33').Value = '  +1.76%  '
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('D35').Value = '150.39'
$ws.Range('E35').Value = '  -3.50%  '
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('D38').Value = '0.896'
$ws.Range('E38').Value = '  +4.92%  '
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').Value = '0.844'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '36.23'
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').Value = '289.94'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('E44').Value = '  +2.13%  '
$ws.Range('D45').Value = '0.621'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0557'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '0.996'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('E51').Value = '  +0.55%  '
